$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A3").Value = 427374227
[void]$ws.Range("A4").Select()
